$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing B2/B3 values
$ws.Range("B2").Value = 1
$ws.Range("B3").Value = 1

# Add new rows 4-7
$data = @(
    @(0, 1, -0.4, -2),
    @(12, 1, -0.4, -2),
    @(0, 1, -1.6, -2),
    @(12, 1, -1.6, -2)
)

$row = 4
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $row++
}

# Update selection to B8 to match diff
$ws.Range("B8").Select()
